$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.045
$ws.Range("A9").Value = -21.831
$ws.Range("C12").Value = -11.623
$ws.Range("E13").Value = 16.717
$ws.Range("D15").Value = -8.348000000000003
$ws.Range("E16").Value = 16.85
$ws.Range("A18").Value = -22.198
$ws.Range("A20").Value = -20.447
$ws.Range("E20").Value = 16.295
$ws.Range("E24").Value = 16.512
$ws.Range("C26").Value = -12.808
$ws.Range("A27").Value = -21.761
$ws.Range("C27").Value = -13.363
$ws.Range("C29").Value = -12.124
$ws.Range("C37").Value = -13.351
$ws.Range("C38").Value = -13.818
$ws.Range("D38").Value = -7.935
$ws.Range("E39").Value = 16.373
$ws.Range("D44").Value = -7.417999999999999
$ws.Range("E48").Value = 17.347
$ws.Range("C51").Value = -12.751
$ws.Range("D51").Value = -7.590000000000001
$ws.Range("E52").Value = 16.759
$ws.Range("C55").Value = -13.65
$ws.Range("E56").Value = 16.881
$ws.Range("D57").Value = -8.032
$ws.Range("D63").Value = -7.337000000000001
$ws.Range("A69").Value = -21.759
$ws.Range("C69").Value = -11.747
$ws.Range("C70").Value = -13.58
$ws.Range("D70").Value = -7.982000000000001
$ws.Range("A76").Value = -20.306
$ws.Range("A82").Value = -22.205
$ws.Range("C83").Value = -13.509
$ws.Range("E84").Value = 16.68
$ws.Range("D99").Value = -8.104000000000001
$ws.Range("E100").Value = 16.741
$ws.Range("E101").Value = 16.667
$ws.Range("C102").Value = -13.419
